$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 8927
$ws.Range("I28").Value = 5559.25
$ws.Range("J28").Value = 10610.875
$ws.Range("K28").Value = 5559.25
$ws.Range("L28").Value = 10610.875
$ws.Range("M28").Value = -5074.25
$ws.Range("N28").Value = -11580.875
$ws.Range("H62").Value = 5714.8887
$ws.Range("I62").Value = 5131.3335
$ws.Range("K62").Value = 5131.3335
$ws.Range("M62").Value = -4507.3335
$ws.Range("H65").Value = 5714.8887
$ws.Range("I65").Value = 5131.3335
$ws.Range("K65").Value = 25656.6675
$ws.Range("M65").Value = -22536.6675
$ws.Range("H96").Value = 726189.25
$ws.Range("J96").Value = 1450968.1
$ws.Range("L96").Value = 4352904.300000001
$ws.Range("N96").Value = -4355650.300000001
$ws.Range("H107").Value = 1382.8889
$ws.Range("I107").Value = 709.2
$ws.Range("K107").Value = 709.2
$ws.Range("M107").Value = 1210.8
$ws.Range("H111").Value = 9871.817999999999
$ws.Range("I111").Value = 4919.2
$ws.Range("J111").Value = 13999
$ws.Range("K111").Value = 14757.6
$ws.Range("L111").Value = 41997
$ws.Range("M111").Value = -11690.6
$ws.Range("N111").Value = -48131
$ws.Range("H112").Value = 1394089.9
$ws.Range("J112").Value = 1476031.1
$ws.Range("L112").Value = 4428093.300000001
$ws.Range("N112").Value = -4430309.300000001
$ws.Range("H125").Value = 2028.75
$ws.Range("I125").Value = 1556.75
$ws.Range("J125").Value = 2500.75
$ws.Range("K125").Value = 14010.75
$ws.Range("L125").Value = 22506.75
$ws.Range("M125").Value = -11550.75
$ws.Range("N125").Value = -27426.75
$ws.Range("H137").Value = 2673.5144
$ws.Range("J137").Value = 4070.1667
$ws.Range("L137").Value = 12210.5001
$ws.Range("N137").Value = -17310.5001
$ws.Range("H141").Value = 13892855
$ws.Range("I141").Value = 15154778
$ws.Range("K141").Value = 45464334
$ws.Range("M141").Value = -45459154
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21985.143
$ws.Range("I28").Value = 8982.833000000001
$ws.Range("K28").Value = 8982.833000000001
$ws.Range("M28").Value = -8790.833000000001
$ws.Range("H32").Value = 3209.6406
$ws.Range("I32").Value = 3107.2678
$ws.Range("J32").Value = 3926.25
$ws.Range("K32").Value = 3107.2678
$ws.Range("L32").Value = 3926.25
$ws.Range("M32").Value = -2820.2678
$ws.Range("N32").Value = -4500.25
$ws.Range("H74").Value = 1506.6666
$ws.Range("I74").Value = 1402.5
$ws.Range("J74").Value = 2340
$ws.Range("K74").Value = 1402.5
$ws.Range("L74").Value = 2340
$ws.Range("M74").Value = -528.5
$ws.Range("N74").Value = -4088
$ws.Range("H77").Value = 1506.6666
$ws.Range("I77").Value = 1402.5
$ws.Range("J77").Value = 2340
$ws.Range("K77").Value = 7012.5
$ws.Range("L77").Value = 11700
$ws.Range("M77").Value = -2644.5
$ws.Range("N77").Value = -20436
$ws.Range("H97").Value = 1558.8422
$ws.Range("I97").Value = 1336.6923
$ws.Range("K97").Value = 1336.6923
$ws.Range("M97").Value = -840.6922999999999
$ws.Range("H99").Value = 21985.143
$ws.Range("I99").Value = 8982.833000000001
$ws.Range("K99").Value = 8982.833000000001
$ws.Range("M99").Value = -5987.833000000001
$ws.Range("H102").Value = 3282.625
$ws.Range("J102").Value = 4462.75
$ws.Range("L102").Value = 4462.75
$ws.Range("N102").Value = -7706.75
$ws.Range("H106").Value = 49999.668
$ws.Range("J106").Value = 49999.668
$ws.Range("L106").Value = 49999.668
$ws.Range("N106").Value = -52523.668
$ws.Range("H110").Value = 4482.7334
$ws.Range("I110").Value = 4264.923
$ws.Range("J110").Value = 5898.5
$ws.Range("K110").Value = 4264.923
$ws.Range("L110").Value = 5898.5
$ws.Range("M110").Value = -2219.923
$ws.Range("N110").Value = -9988.5
$ws.Range("H112").Value = 114999.5
$ws.Range("J112").Value = 114999.5
$ws.Range("L112").Value = 114999.5
$ws.Range("N112").Value = -117953.5
$ws.Range("H122").Value = 2891.7292
$ws.Range("I122").Value = 2730.5676
$ws.Range("J122").Value = 3433.818
$ws.Range("K122").Value = 8191.702799999999
$ws.Range("L122").Value = 10301.454
$ws.Range("M122").Value = -5741.702799999999
$ws.Range("N122").Value = -15201.454
$ws.Range("H124").Value = 69000
$ws.Range("J124").Value = 69000
$ws.Range("L124").Value = 69000
$ws.Range("N124").Value = -78820
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1117.25
$ws.Range("I64").Value = 984.5
$ws.Range("K64").Value = 984.5
$ws.Range("M64").Value = -759.5
$ws.Range("H67").Value = 1117.25
$ws.Range("I67").Value = 984.5
$ws.Range("K67").Value = 984.5
$ws.Range("M67").Value = -204.5
$ws.Range("H86").Value = 5929.5557
$ws.Range("I86").Value = 1796
$ws.Range("K86").Value = 1796
$ws.Range("M86").Value = -673
$ws.Range("H89").Value = 5929.5557
$ws.Range("I89").Value = 1796
$ws.Range("K89").Value = 8980
$ws.Range("M89").Value = -3364
$ws.Range("H94").Value = 1961.6428
$ws.Range("I94").Value = 2064.25
$ws.Range("K94").Value = 2064.25
$ws.Range("M94").Value = -1613.25
$ws.Range("H97").Value = 25027850
$ws.Range("I97").Value = 33337134
$ws.Range("K97").Value = 33337134
$ws.Range("M97").Value = -33336143
$ws.Range("H105").Value = 718046
$ws.Range("I105").Value = 1431451.5
$ws.Range("J105").Value = 4640.5
$ws.Range("K105").Value = 1431451.5
$ws.Range("L105").Value = 4640.5
$ws.Range("M105").Value = -1429704.5
$ws.Range("N105").Value = -8134.5
$ws.Range("H134").Value = 4349590.5
$ws.Range("I134").Value = 1856.238
$ws.Range("J134").Value = 50000800
$ws.Range("K134").Value = 5568.714
$ws.Range("L134").Value = 150002400
$ws.Range("M134").Value = -3033.714
$ws.Range("N134").Value = -150007470
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5189.032
$ws.Range("I16").Value = 5867.3076
$ws.Range("K16").Value = 5867.3076
$ws.Range("M16").Value = -5580.3076
$ws.Range("H31").Value = 19610390
$ws.Range("I31").Value = 31252688
$ws.Range("J31").Value = 2309.842
$ws.Range("K31").Value = 31252688
$ws.Range("L31").Value = 2309.842
$ws.Range("M31").Value = -31252393
$ws.Range("N31").Value = -2899.842
$ws.Range("H34").Value = 19610390
$ws.Range("I34").Value = 31252688
$ws.Range("J34").Value = 2309.842
$ws.Range("K34").Value = 31252688
$ws.Range("L34").Value = 2309.842
$ws.Range("M34").Value = -31252486
$ws.Range("N34").Value = -2713.842
$ws.Range("H58").Value = 3732.6365
$ws.Range("I58").Value = 3514.5
$ws.Range("K58").Value = 3514.5
$ws.Range("M58").Value = -3311.5
$ws.Range("H113").Value = 5189.032
$ws.Range("I113").Value = 5867.3076
$ws.Range("K113").Value = 5867.3076
$ws.Range("M113").Value = -3697.3076
$ws.Range("H122").Value = 4609.2666
$ws.Range("I122").Value = 4275.5
$ws.Range("J122").Value = 4990.7144
$ws.Range("K122").Value = 12826.5
$ws.Range("L122").Value = 14972.1432
$ws.Range("M122").Value = -10376.5
$ws.Range("N122").Value = -19872.1432
$ws.Range("H136").Value = 3732.6365
$ws.Range("I136").Value = 3514.5
$ws.Range("K136").Value = 10543.5
$ws.Range("M136").Value = -7993.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 734
$ws.Range("I86").Value = 601.6667
$ws.Range("J86").Value = 866.3333
$ws.Range("K86").Value = 1805.0001
$ws.Range("L86").Value = 2598.9999
$ws.Range("M86").Value = -619.0001
$ws.Range("N86").Value = -4970.9999
$ws.Range("H89").Value = 734
$ws.Range("I89").Value = 601.6667
$ws.Range("J89").Value = 866.3333
$ws.Range("K89").Value = 5415.0003
$ws.Range("L89").Value = 7796.9997
$ws.Range("M89").Value = 512.9997000000003
$ws.Range("N89").Value = -19652.9997
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 268
$ws.Range("I2").Value = 349
$ws.Range("K2").Value = 349
$ws.Range("M2").Value = -236
$ws.Range("H113").Value = 2061476
$ws.Range("I113").Value = 3996
$ws.Range("J113").Value = 3090216
$ws.Range("K113").Value = 3996
$ws.Range("L113").Value = 3090216
$ws.Range("M113").Value = -1826
$ws.Range("N113").Value = -3094556
$ws.Range("H132").Value = 27276012
$ws.Range("I132").Value = 4472.25
$ws.Range("K132").Value = 13416.75
$ws.Range("M132").Value = -10886.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6791.625
$ws.Range("I7").Value = 5722.1665
$ws.Range("K7").Value = 5722.1665
$ws.Range("M7").Value = -5610.1665
$ws.Range("H93").Value = 2782763
$ws.Range("I93").Value = 2337.3635
$ws.Range("K93").Value = 2337.3635
$ws.Range("M93").Value = -1089.3635
$ws.Range("H100").Value = 10013006
$ws.Range("I100").Value = 2728.6428
$ws.Range("K100").Value = 2728.6428
$ws.Range("M100").Value = -2187.6428
$ws.Range("H126").Value = 6791.625
$ws.Range("I126").Value = 5722.1665
$ws.Range("K126").Value = 17166.4995
$ws.Range("M126").Value = -14696.4995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5406
$ws.Range("H81").Value = 1445.6
$ws.Range("I81").Value = 1445.6
$ws.Range("K81").Value = 2891.2
$ws.Range("M81").Value = -1830.2
$ws.Range("H84").Value = 1445.6
$ws.Range("I84").Value = 1445.6
$ws.Range("K84").Value = 14456
$ws.Range("M84").Value = -9152
$ws.Range("H126").Value = 2760.0312
$ws.Range("I126").Value = 2843.3076
$ws.Range("J126").Value = 2399.1667
$ws.Range("K126").Value = 8529.9228
$ws.Range("L126").Value = 7197.500100000001
$ws.Range("M126").Value = -6059.9228
$ws.Range("N126").Value = -12137.5001
